$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Stations_VenueCounts")

# --- Row 6 (Lees) ---
$ws.Range("J6").Value = 14
$ws.Range("M6").Value = 106
$ws.Range("N6").Value = 61

# --- Row 10 (Rideau) ---
$ws.Range("E10").Value = 0
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 13
$ws.Range("H10").Value = 2
$ws.Range("I10").Value = 7
$ws.Range("J10").Value = 41
$ws.Range("K10").Value = 15
$ws.Range("L10").Value = 5
$ws.Range("M10").Value = 90
$ws.Range("N10").Value = 50

# --- Row 13 (uOttawa) ---
$ws.Range("E13").Value = 1
$ws.Range("F13").Value = 83
$ws.Range("G13").Value = 11
$ws.Range("H13").Value = 4
$ws.Range("I13").Value = 6
$ws.Range("J13").Value = 8
$ws.Range("K13").Value = 5
$ws.Range("L13").Value = 4
$ws.Range("M13").Value = 131
$ws.Range("N13").Value = 50

# --- Row 16 (Dundas) ---
$ws.Range("K16").Value = 47
$ws.Range("N16").Value = 99

# --- Row 18 (Museum) ---
$ws.Range("I18").Value = 12
$ws.Range("M18").Value = 164
$ws.Range("N18").Value = 82

# --- Static (non-formula) "Unique Categories" totals that don't auto-recalc ---
$ws.Range("N32").Value = 290
$ws.Range("N43").Value = 393

# --- View state: scroll the window so G20 is the top-left visible cell ---
# (selection/activeCell itself - N33 - is left untouched)
$win = $excel.ActiveWindow
$win.ScrollRow = 20
$win.ScrollColumn = 7
$win.TopLeftCell = $ws.Range("G20")

$excel.Calculate()
